$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (shifts existing rows 12-30 down to 13-31)
$row = $ws.Rows.Item(12)
$row.Insert()

# Populate the new row 12 with the new weekly data point
$ws.Cells.Item(12,1).Value  = 11
$ws.Cells.Item(12,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(12,3).Value  = "Bíobío"
$ws.Cells.Item(12,4).Value  = 44649
$ws.Cells.Item(12,5).Value  = 8
$ws.Cells.Item(12,6).Value  = 100112037
$ws.Cells.Item(12,7).Value  = "Cebollín"
$ws.Cells.Item(12,8).Value  = "Sin especificar"
$ws.Cells.Item(12,9).Value  = "Primera"
$ws.Cells.Item(12,10).Value = 220
$ws.Cells.Item(12,11).Value = 8000
$ws.Cells.Item(12,12).Value = 8500
$ws.Cells.Item(12,13).Value = 8227
$ws.Cells.Item(12,14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(12,15).Value = "Región Metropolitana"
$ws.Cells.Item(12,16).Value = 229
$ws.Cells.Item(12,17).Value = 36
$ws.Cells.Item(12,18).Value = "Hortaliza"
